# Updated capital structure database
# Applies the revised ratio figures to the two data rows (2 and 3) of the
# "earnings_debt" sheet: operating/net margins, cash-return metrics, ROE/ROIC
# spreads and net-debt ratios, and drops the now-unused "buybacks_cash_returned"
# (column T) value while adding the new roe_cost_equity / sales_invested_capital
# / roic / roic_cost_capital figures (columns W, Y, Z, AA, AC).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = [ordered]@{
    "I" = 0.676470588235294
    "J" = 0.676470588235294
    "K" = 0.175
    "L" = 0.6433823529411764
    "M" = 0
    "N" = 0
    "O" = 0
    "P" = 0
    "Q" = 0
    "R" = 0
    "U" = 0.59
    "V" = 0.025
    "W" = 0.01762336354481369
    "X" = 0.09238291079777955
    "Y" = -0.07475954725296585
    "Z" = 0.02777777777777778
    "AA" = 0.01879084967320261
    "AB" = 0.09238291079777955
    "AC" = -0.07359206112457693
    "AG" = -0.59
    "AJ" = -0.02564102564102564
    "AK" = -0.03087388801674516
}

foreach ($row in 2, 3) {
    # Column T (buybacks_cash_returned) no longer has a value for these rows.
    $ws.Range("T$row").ClearContents()

    foreach ($col in $newValues.Keys) {
        $ws.Range("$col$row").Value = $newValues[$col]
    }
}
